# Fruta / hortaliza, semanal
# A new weekly record is inserted at row 306 (pushing the existing
# rows 306-359 down to 307-360), with the new record's own values set
# afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 306; Excel shifts rows 306:359 down to 307:360,
# copying values and formatting (matches the observed row-shift diff).
$ws.Rows("306:306").Insert()

# Populate the newly inserted row with this week's data.
$ws.Range("A306").Value2 = 10
$ws.Range("B306").Value2 = "Vega Modelo de Temuco"
$ws.Range("C306").Value2 = "La Araucanía"
$ws.Range("D306").Value2 = 45258
$ws.Range("E306").Value2 = 9
$ws.Range("F306").Value2 = 100112005
$ws.Range("G306").Value2 = "Puerro"
$ws.Range("H306").Value2 = "Azul de Maquehue"
$ws.Range("I306").Value2 = "Primera"
$ws.Range("J306").Value2 = 110
$ws.Range("K306").Value2 = 12000
$ws.Range("L306").Value2 = 12000
$ws.Range("M306").Value2 = 12000
$ws.Range("N306").Value2 = "`$/docena de paquetes"
$ws.Range("O306").Value2 = "Provincia de Cautín"
$ws.Range("P306").Value2 = 1000
$ws.Range("Q306").Value2 = 12
$ws.Range("R306").Value2 = "Hortaliza"
